$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "21.544.21"
$ws.Range("E2").Value = "  -2.49%  "

$ws.Range("D3").Value = "1.530.83"
$ws.Range("E3").Value = "  -1.74%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'1.002"
$ws.Range("E5").Value = "  +0.22%  "

$ws.Range("D6").Value = "'288.76"
$ws.Range("E6").Value = "  -1.13%  "

$ws.Range("D7").Value = "'0.3868"
$ws.Range("E7").Value = "  -2.72%  "

$ws.Range("D8").Value = "'0.3169"
$ws.Range("E8").Value = "  -2.18%  "

$ws.Range("D9").Value = "'42.73"
$ws.Range("E9").Value = "  -3.36%  "

$ws.Range("D10").Value = "'0.07152"
$ws.Range("E10").Value = "  -2.45%  "

$ws.Range("D11").Value = "'1.065"
$ws.Range("E11").Value = "  -2.07%  "

$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "  +0.16%  "

$ws.Range("D13").Value = "'5.729"
$ws.Range("E13").Value = "  +0.76%  "

$ws.Range("D14").Value = "'18.16"
$ws.Range("E14").Value = "  -4.61%  "

$ws.Range("D15").Value = "'6.550"
$ws.Range("E15").Value = "  -1.67%  "

$ws.Range("D16").Value = "1.539.63"
$ws.Range("E16").Value = "  -1.11%  "

$ws.Range("D17").Value = "'0.00001087"
$ws.Range("E17").Value = "  -5.21%  "

$ws.Range("D18").Value = "'0.06617"
$ws.Range("E18").Value = "  +0.17%  "

$ws.Range("D19").Value = "'83.58"
$ws.Range("E19").Value = "  -0.35%  "

$ws.Range("E20").Value = "  +0.27%  "

$ws.Range("D21").Value = "'6.097"
$ws.Range("E21").Value = "  -3.56%  "

$ws.Range("D22").Value = "'15.39"
$ws.Range("E22").Value = "  -2.53%  "

$ws.Range("D23").Value = "'10.80"
$ws.Range("E23").Value = "  -4.46%  "

$ws.Range("D24").Value = "'2.371"
$ws.Range("E24").Value = "  +1.20%  "

$ws.Range("D25").Value = "21.543.61"
$ws.Range("E25").Value = "  -2.56%  "

$ws.Range("D26").Value = "'2.379"
$ws.Range("E26").Value = "  -2.64%  "

$ws.Range("D27").Value = "'149.84"
$ws.Range("E27").Value = "  +0.73%  "

$ws.Range("E28").Value = "  -1.58%  "

$ws.Range("D29").Value = "'4.834"
$ws.Range("E29").Value = "  -0.90%  "

$ws.Range("D30").Value = "1.708.57"
$ws.Range("E30").Value = "  -1.34%  "

$ws.Range("D31").Value = "'116.51"
$ws.Range("E31").Value = "  -2.19%  "

$ws.Range("D32").Value = "'6.040"
$ws.Range("E32").Value = "  +5.18%  "

$ws.Range("D33").Value = "'0.9492"
$ws.Range("E33").Value = "  -5.93%  "

$ws.Range("D34").Value = "'0.07995"
$ws.Range("E34").Value = "  -4.57%  "

$ws.Range("D35").Value = "'8.522"
$ws.Range("E35").Value = "  -5.93%  "

$ws.Range("D36").Value = "'5.157"
$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("D37").Value = "'1.485"
$ws.Range("E37").Value = "  -8.75%  "

$ws.Range("D38").Value = "'0.02206"
$ws.Range("E38").Value = "  -3.48%  "

$ws.Range("D39").Value = "'11.27"
$ws.Range("E39").Value = "  +4.46%  "

$ws.Range("D40").Value = "'0.05886"
$ws.Range("E40").Value = "  -4.36%  "

$ws.Range("D41").Value = "'0.2022"
$ws.Range("E41").Value = "  -2.33%  "

$ws.Range("D42").Value = "'1.178"
$ws.Range("E42").Value = "  -3.33%  "

$ws.Range("E43").Value = "  +0.27%  "

$ws.Range("D44").Value = "'0.5744"
$ws.Range("E44").Value = "  -2.30%  "

$ws.Range("D45").Value = "'13.20"
$ws.Range("E45").Value = "  +0.69%  "

$ws.Range("D46").Value = "'3.719"
$ws.Range("E46").Value = "  -1.22%  "

$ws.Range("D47").Value = "'0.5546"
$ws.Range("E47").Value = "  -1.45%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.888"
$ws.Range("E48").Value = "  -0.98%  "

$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "'1.157"
$ws.Range("E49").Value = "  +1.36%  "

$ws.Range("D50").Value = "'115.42"
$ws.Range("E50").Value = "  -3.16%  "

$ws.Range("D51").Value = "'0.06680"
$ws.Range("E51").Value = "  -2.62%  "
